$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string / row of data (row 5)
$ws.Range("A5").Value = 44317
$ws.Range("A5").NumberFormat = "d-mmm"

$ws.Range("B5").Value = 0.052083333333333336
$ws.Range("B5").NumberFormat = "h:mm"

$ws.Range("C5").Value = "apprence - travail sur les avatars"

# Update selection to match the saved state (C6)
$ws.Range("C6").Select()
